$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.06992366666666666
$ws.Range("H2").Value = 0.209771
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.06247066666666667
$ws.Range("N2").Value = 0.187412
$ws.Range("O2").Value = 0.009697282539039283
$ws.Range("P2").Value = 0.009697282539039282
$ws.Range("Q2").Value = 0.004368178072444444
$ws.Range("R2").Value = 0.039313602652
$ws.Range("S2").Value = 0.009697282539039283
$ws.Range("T2").Value = 0.009697282539039282

$ws.Range("G3").Value = 0.06992366666666666
$ws.Range("H3").Value = 0.209771
$ws.Range("O3").Value = 0.5253497589468907
$ws.Range("P3").Value = 0.5253497589468907
$ws.Range("Q3").Value = 0.236645811665
$ws.Range("R3").Value = 2.129812304985
$ws.Range("S3").Value = 0.5253497589468907
$ws.Range("T3").Value = 0.5253497589468907

$ws.Range("G4").Value = 0.06992366666666666
$ws.Range("H4").Value = 0.209771
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.995264
$ws.Range("N4").Value = 8.985792
$ws.Range("O4").Value = 0.46495295851407
$ws.Range("P4").Value = 0.4649529585140699
$ws.Range("Q4").Value = 0.2094398415146667
$ws.Range("R4").Value = 1.884958573632
$ws.Range("S4").Value = 0.46495295851407
$ws.Range("T4").Value = 0.4649529585140699
